$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.532.53'

$ws.Range("D3").Value = '1.846.72'
$ws.Range("E3").Value = '  -4.20%  '

$ws.Range("E4").Value = '  -1.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.39'
$ws.Range("E5").Value = '  +2.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4644'
$ws.Range("E7").Value = '  -4.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3892'
$ws.Range("E8").Value = '  -4.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.19'
$ws.Range("E9").Value = '  -2.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07881'
$ws.Range("E10").Value = '  -3.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9738'
$ws.Range("E11").Value = '  -3.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.15'
$ws.Range("E12").Value = '  -6.97%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.807'
$ws.Range("E13").Value = '  -4.60%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.814.91'
$ws.Range("E14").Value = '  -6.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.958'
$ws.Range("E15").Value = '  -4.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06927'
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("E17").Value = '  -0.97%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.69'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001002'
$ws.Range("E19").Value = '  -3.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.02'
$ws.Range("E20").Value = '  -3.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  -0.73%  '

$ws.Range("D22").Value = '28.571.10'
$ws.Range("E22").Value = '  -3.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.363'
$ws.Range("E23").Value = '  -5.40%  '

$ws.Range("E24").Value = '  -7.12%  '

$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("D26").Value = '2.101.59'
$ws.Range("E26").Value = '  -2.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.41'
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.38'
$ws.Range("E28").Value = '  -3.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.951'
$ws.Range("E29").Value = '  -7.89%  '

$ws.Range("E30").Value = '  -4.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.50'
$ws.Range("E31").Value = '  -2.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9632'
$ws.Range("E32").Value = '  -5.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09373'
$ws.Range("E33").Value = '  -2.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.352'
$ws.Range("E34").Value = '  -4.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.456'
$ws.Range("E35").Value = '  -2.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.341'
$ws.Range("E36").Value = '  -4.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06065'
$ws.Range("E37").Value = '  -6.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02187'
$ws.Range("E38").Value = '  -4.31%  '

$ws.Range("E39").Value = '  -3.83%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5679'
$ws.Range("E40").Value = '  -4.43%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.636'
$ws.Range("E41").Value = '  -3.38%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.09'
$ws.Range("E42").Value = '  -6.13%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1792'
$ws.Range("E43").Value = '  -2.93%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.382'
$ws.Range("E44").Value = '  -6.34%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.248'
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5351'
$ws.Range("E46").Value = '  -3.84%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.72'
$ws.Range("E47").Value = '  -5.59%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07067'
$ws.Range("E48").Value = '  -6.01%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.880'
$ws.Range("E49").Value = '  -4.37%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.24'
$ws.Range("E50").Value = '  -4.51%  '

$ws.Range("B51").Value = 'Chiliz'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GSCt2y6YSgO26+chiliz-chz'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1296'
$ws.Range("E51").Value = '  -0.22%  '
